$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168023824691772
$ws.Range("B1").Value = 2.375645160675049
$ws.Range("D1").Value = 2.385340213775635
$ws.Range("E1").Value = 1.212505340576172
